# Update the simulation output data (columns B:K, rows 2:11) on Sheet1 with
# the re-run values ("modifications made to simulation constraints"), and
# narrow the stray full-column selection from M1:X1048576 to M1:W1048576.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows 2-11, columns B(2)-K(11)
$data = @(
    @(865,384,224,144,96,64,41,24,10,0),
    @(864,384,224,144,96,64,41,24,11,11),
    @(864,384,224,144,96,64,41,24,24,24),
    @(864,384,224,144,96,64,41,42,42,42),
    @(864,384,224,144,96,64,65,66,66,66),
    @(864,384,224,144,96,97,98,99,99,99),
    @(864,384,224,144,145,146,147,149,149,149),
    @(848,374,216,216,216,216,216,216,216,216),
    @(532,216,216,216,216,216,216,216,216,216),
    @(216,212,212,212,212,212,212,212,212,212)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $colNum = $j + 2
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$j]
    }
}

# Shrink the selected "rest of the sheet" range by one column (was M1:X1048576).
$ws.Range("M1:W1048576").Select()
